# Fill in timesheet hours for the week of 15-19 Jan (columns L,O,P,Q,R)
# across the task rows (9-33), matching the commit's "Add files via upload".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 - Unit Test Plan Preparation
$ws.Range("L14").Value = 1
$ws.Range("O14").Value = 1

# Row 19 - LLD Rework
$ws.Range("O19").Value = 1.5

# Row 24 - Code Peer Review
$ws.Range("Q24").Value = 1

# Row 25 - Peer Testing
$ws.Range("Q25").Value = 1

# Row 26 - Test Result review
$ws.Range("Q26").Value = 1

# Row 27 - Rework
$ws.Range("Q27").Value = 2

# Row 28 - Code Integration
$ws.Range("Q28").Value = 2

# Row 29 - Integration Testing
$ws.Range("Q29").Value = 1

# Row 30 - Rework
$ws.Range("Q30").Value = 1

# Row 31 - Test Result review
$ws.Range("Q31").Value = 1

# Row 32 - Presentation Preparation
$ws.Range("R32").Value = 3

# Row 33 - Deployment
$ws.Range("R33").Value = 1

# Move the selection/cursor as it ended up after the edits
$ws.Range("G15").Select()
